$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45202 = 2023-10-03)
# that was bumped by one day (45203 = 2023-10-04) for every data row
# (rows 2 through 353).
$ws.Range("C2:C353").Value = 45203
